# This workbook holds a weekly time series of Cebollín prices for
# "Feria Lagunitas de Puerto Montt". A new weekly record needs to be
# inserted as the new row 338, which pushes every existing record
# starting at the old row 338 down by one row (old row 375 becomes the
# new row 376).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at 338; this shifts rows 338:375 down to
# 339:376 and grows the used range to A1:R376 automatically.
$ws.Rows("338:338").Insert()

# Populate the newly inserted row 338 with the new weekly record.
# Columns A, B, C, E, F, G, H, I, N, O, Q, R keep the same
# market/category metadata as the surrounding rows; D (date), J
# (volume), K/L/M (min/max/avg price) and P (price per kg) hold the
# new observation.
$ws.Range("A338").Value = 4
$ws.Range("B338").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C338").Value = "Los Lagos"
$ws.Range("D338").NumberFormat = $ws.Range("D339").NumberFormat
$ws.Range("D338").Value = 44918
$ws.Range("E338").Value = 10
$ws.Range("F338").Value = 100112037
$ws.Range("G338").Value = "Cebollín"
$ws.Range("H338").Value = "Sin especificar"
$ws.Range("I338").Value = "Primera"
$ws.Range("J338").Value = 100
$ws.Range("K338").Value = 7000
$ws.Range("L338").Value = 7000
$ws.Range("M338").Value = 7000
$ws.Range("N338").Value = "`$/paquete 36 unidades"
$ws.Range("O338").Value = "Región Metropolitana"
$ws.Range("P338").Value = 194
$ws.Range("Q338").Value = 36
$ws.Range("R338").Value = "Hortaliza"
